# Solucionario_parcial.xlsx edit
#
# Renames the original sheet to "Parcial 1" and adds a second exam variant
# "Parcial 1 _ version 2" (duplicated from the first so all of the original
# bold / fill / border formatting carries over), then edits the duplicate's
# answers and truth-table content to match version 2 of the exam.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Parcial 1"

# Duplicate "Parcial 1" right after itself so every style (bold answer
# column, shaded question column, bordered truth-table cells) is carried
# over exactly as authored.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Parcial 1 _ version 2"

# ---------------------------------------------------------------------
# Rebuild the two truth tables on the duplicate sheet.
# ---------------------------------------------------------------------
# Slide the existing second table (H10:M14 -> p,q,~p,~q,(~pv~q),(~pv~q)^q)
# two columns right, to J10:O14, keeping all of its formatting.
$ws2.Range("H10:M14").Copy($ws2.Range("J10"))

# Column I must end up completely unused in the new layout -- clear
# formats before contents so no placeholder cell is left behind.
$ws2.Range("I10:I14").ClearFormats()
$ws2.Range("I10:I14").ClearContents()

# Build the brand-new first table in columns G:H using the same
# header/data style as the existing C:D columns.
$ws2.Range("C10:D14").Copy($ws2.Range("G10"))

# ---------------------------------------------------------------------
# Answers column (B1:B10) for exam version 2.
# ---------------------------------------------------------------------
$ws2.Range("B1").Value2  = "Diagrama de flujo"
$ws2.Range("B2").Value2  = "B"
$ws2.Range("B3").Value2  = "B"
$ws2.Range("B4").Value2  = "A"
$ws2.Range("B5").Value2  = "C"
$ws2.Range("B6").Value2  = "B"
$ws2.Range("B7").Value2  = "A"
$ws2.Range("B8").Value2  = "D"
$ws2.Range("B9").Value2  = "A"
$ws2.Range("B10").Value2 = "A"

# B11:B14 stay empty but keep the bold formatting of the rest of column B.
$ws2.Range("B11:B14").Font.Bold = $true

# ---------------------------------------------------------------------
# Truth-table headers (row 10).
# ---------------------------------------------------------------------
$ws2.Range("C10").Value2 = "p"
$ws2.Range("D10").Value2 = "q"
$ws2.Range("E10").Value2 = "~p"
$ws2.Range("F10").Value2 = "~q"
$ws2.Range("G10").Value2 = "(p^¬q)"
$ws2.Range("H10").Value2 = "(p^¬q)v¬p"

$ws2.Range("J10").Value2 = "p"
$ws2.Range("K10").Value2 = "q"
$ws2.Range("L10").Value2 = "~p"
$ws2.Range("M10").Value2 = "~q"
$ws2.Range("N10").Value2 = "(~pv~q)"
$ws2.Range("O10").Value2 = "(¬pv¬q)^q"

# ---------------------------------------------------------------------
# Truth-table values (rows 11:14).
# ---------------------------------------------------------------------
$tableLeft = @(
  @("v","v","f","f","f","f"),
  @("v","f","f","v","v","v"),
  @("f","v","v","f","f","v"),
  @("f","f","v","v","f","v")
)
$tableRight = @(
  @("v","v","f","f","f","f"),
  @("v","f","f","v","v","f"),
  @("f","v","v","f","v","v"),
  @("f","f","v","v","v","f")
)
$leftCols  = @("C","D","E","F","G","H")
$rightCols = @("J","K","L","M","N","O")

for ($r = 0; $r -lt 4; $r++) {
  $row = 11 + $r
  for ($c = 0; $c -lt 6; $c++) {
    $ws2.Range("$($leftCols[$c])$row").Value2  = $tableLeft[$r][$c]
    $ws2.Range("$($rightCols[$c])$row").Value2 = $tableRight[$r][$c]
  }
}

# ---------------------------------------------------------------------
# Column widths: new layout custom-widths columns A, B, H and O.
# ---------------------------------------------------------------------
$oldMWidth = $ws2.Columns.Item(13).ColumnWidth
$ws2.Columns.Item(2).ClearFormats()
$ws2.Range("B1:B14").Font.Bold = $true
$ws2.Columns.Item(8).ColumnWidth  = $oldMWidth
$ws2.Columns.Item(15).ColumnWidth = $oldMWidth
$ws2.Columns.Item(13).ClearFormats()

# ---------------------------------------------------------------------
# Selections / active sheet, matching the authored state.
# ---------------------------------------------------------------------
$ws1.Range("B23").Select() | Out-Null
$ws2.Range("O12").Select() | Out-Null

Write-Host "Sheets now:"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
  Write-Host "  $i -> $($wb.Worksheets.Item($i).Name)"
}
